$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the greeting text for rule R10 (cell E8) from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Make the sheet active and select E8, matching the saved selection state
$ws.Activate()
$ws.Range("E8").Select()
